$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 74453.63
$ws.Range("J123").Value = 74453.63
$ws.Range("L123").Value = 74453.63
$ws.Range("N123").Value = -84253.63
$ws.Range("H138").Value = 2139.9883
$ws.Range("I138").Value = 1177.3715
$ws.Range("J138").Value = 2813.82
$ws.Range("K138").Value = 3532.1145
$ws.Range("L138").Value = 8441.460000000001
$ws.Range("M138").Value = 1607.8855
$ws.Range("N138").Value = -18721.46

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4809.17
$ws.Range("I32").Value = 4035.821
$ws.Range("J32").Value = 19502.8
$ws.Range("K32").Value = 4035.821
$ws.Range("L32").Value = 19502.8
$ws.Range("M32").Value = -3748.821
$ws.Range("N32").Value = -20076.8
$ws.Range("H44").Value = 18450
$ws.Range("I44").Value = 2400
$ws.Range("J44").Value = 24870
$ws.Range("K44").Value = 2400
$ws.Range("L44").Value = 24870
$ws.Range("M44").Value = -1912
$ws.Range("N44").Value = -25846
$ws.Range("H54").Value = 18100
$ws.Range("J54").Value = 18100
$ws.Range("L54").Value = 18100
$ws.Range("N54").Value = -19638
$ws.Range("H61").Value = 291181.16
$ws.Range("I61").Value = 5904.5356
$ws.Range("K61").Value = 5904.5356
$ws.Range("M61").Value = -5692.5356
$ws.Range("H132").Value = 2086231.2
$ws.Range("I132").Value = 2091.647
$ws.Range("K132").Value = 6274.941
$ws.Range("M132").Value = -3744.941
$ws.Range("H134").Value = 56686
$ws.Range("J134").Value = 56686
$ws.Range("L134").Value = 56686
$ws.Range("N134").Value = -66826
$ws.Range("H136").Value = 291181.16
$ws.Range("I136").Value = 5904.5356
$ws.Range("K136").Value = 17713.6068
$ws.Range("M136").Value = -15163.6068
$ws.Range("H137").Value = 35780
$ws.Range("J137").Value = 35780
$ws.Range("L137").Value = 35780
$ws.Range("N137").Value = -45980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 39995.883
$ws.Range("J133").Value = 39995.883
$ws.Range("L133").Value = 39995.883
$ws.Range("N133").Value = -50115.883
$ws.Range("H134").Value = 39156.266
$ws.Range("I134").Value = 8268.294
$ws.Range("J134").Value = 79548.234
$ws.Range("K134").Value = 24804.882
$ws.Range("L134").Value = 238644.702
$ws.Range("M134").Value = -22269.882
$ws.Range("N134").Value = -243714.702

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 57401
$ws.Range("J4").Value = 49251
$ws.Range("L4").Value = 49251
$ws.Range("N4").Value = -49475
$ws.Range("H31").Value = 4493.2256
$ws.Range("I31").Value = 2047.561
$ws.Range("J31").Value = 9268.096
$ws.Range("K31").Value = 2047.561
$ws.Range("L31").Value = 9268.096
$ws.Range("M31").Value = -1752.561
$ws.Range("N31").Value = -9858.096
$ws.Range("H34").Value = 4493.2256
$ws.Range("I34").Value = 2047.561
$ws.Range("J34").Value = 9268.096
$ws.Range("K34").Value = 2047.561
$ws.Range("L34").Value = 9268.096
$ws.Range("M34").Value = -1845.561
$ws.Range("N34").Value = -9672.096
$ws.Range("H58").Value = 210170.89
$ws.Range("I58").Value = 1533.2632
$ws.Range("J58").Value = 346864.53
$ws.Range("K58").Value = 1533.2632
$ws.Range("L58").Value = 346864.53
$ws.Range("M58").Value = -1330.2632
$ws.Range("N58").Value = -347270.53
$ws.Range("H132").Value = 3286.25
$ws.Range("I132").Value = 1470.6666
$ws.Range("J132").Value = 4375.6
$ws.Range("K132").Value = 4411.9998
$ws.Range("L132").Value = 13126.8
$ws.Range("M132").Value = -1881.9998
$ws.Range("N132").Value = -18186.8
$ws.Range("H133").Value = 27659.334
$ws.Range("J133").Value = 27659.334
$ws.Range("L133").Value = 27659.334
$ws.Range("N133").Value = -32719.334
$ws.Range("H134").Value = 185323.81
$ws.Range("I134").Value = 3740.6858
$ws.Range("J134").Value = 503094.3
$ws.Range("K134").Value = 11222.0574
$ws.Range("L134").Value = 1509282.9
$ws.Range("M134").Value = -8687.057400000002
$ws.Range("N134").Value = -1514352.9
$ws.Range("H135").Value = 30503.95
$ws.Range("I135").Value = 30709
$ws.Range("J135").Value = 30493.158
$ws.Range("K135").Value = 30709
$ws.Range("L135").Value = 30493.158
$ws.Range("M135").Value = -25639
$ws.Range("N135").Value = -40633.158
$ws.Range("H136").Value = 210170.89
$ws.Range("I136").Value = 1533.2632
$ws.Range("J136").Value = 346864.53
$ws.Range("K136").Value = 4599.7896
$ws.Range("L136").Value = 1040593.59
$ws.Range("M136").Value = -2049.7896
$ws.Range("N136").Value = -1045693.59
$ws.Range("H140").Value = 37756.25
$ws.Range("I140").Value = 4000
$ws.Range("J140").Value = 38845.16
$ws.Range("K140").Value = 4000
$ws.Range("L140").Value = 38845.16
$ws.Range("M140").Value = 1180
$ws.Range("N140").Value = -49205.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 376285.75
$ws.Range("I4").Value = 750171.7
$ws.Range("J4").Value = 2399.8333
$ws.Range("K4").Value = 2250515.1
$ws.Range("L4").Value = 7199.499899999999
$ws.Range("M4").Value = -2250403.1
$ws.Range("N4").Value = -7423.499899999999
$ws.Range("H131").Value = 1786736.8
$ws.Range("J131").Value = 1075.7234
$ws.Range("L131").Value = 3227.1702
$ws.Range("N131").Value = -13307.1702

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8500
$ws.Range("J5").Value = 8500
$ws.Range("L5").Value = 8500
$ws.Range("N5").Value = -8724
$ws.Range("H122").Value = 31321024
$ws.Range("I122").Value = 50707492
$ws.Range("J122").Value = 4423.385
$ws.Range("K122").Value = 152122476
$ws.Range("L122").Value = 13270.155
$ws.Range("M122").Value = -152120026
$ws.Range("N122").Value = -18170.155
$ws.Range("H132").Value = 6149.657
$ws.Range("I132").Value = 7164.1665
$ws.Range("J132").Value = 3936.182
$ws.Range("K132").Value = 21492.4995
$ws.Range("L132").Value = 11808.546
$ws.Range("M132").Value = -18962.4995
$ws.Range("N132").Value = -16868.546
$ws.Range("H135").Value = 45750.305
$ws.Range("J135").Value = 45750.305
$ws.Range("L135").Value = 45750.305
$ws.Range("N135").Value = -55890.305

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4768691.5
$ws.Range("I122").Value = 5501951.5
$ws.Range("K122").Value = 16505854.5
$ws.Range("M122").Value = -16503404.5
$ws.Range("H132").Value = 17552932
$ws.Range("I132").Value = 23819734
$ws.Range("K132").Value = 71459202
$ws.Range("M132").Value = -71456672
$ws.Range("H134").Value = 28507.908
$ws.Range("J134").Value = 28507.908
$ws.Range("L134").Value = 28507.908
$ws.Range("N134").Value = -38647.908
$ws.Range("H136").Value = 8833.3125
$ws.Range("I136").Value = 8322.611000000001
$ws.Range("J136").Value = 9489.929
$ws.Range("K136").Value = 24967.833
$ws.Range("L136").Value = 28469.787
$ws.Range("M136").Value = -22417.833
$ws.Range("N136").Value = -33569.787
$ws.Range("H139").Value = 41687.855
$ws.Range("J139").Value = 41687.855
$ws.Range("L139").Value = 41687.855
$ws.Range("N139").Value = -51967.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 21733.834
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 21733.834
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 21733.834
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -21957.834
$ws.Range("H122").Value = 1778.3334
$ws.Range("I122").Value = 1712.6923
$ws.Range("K122").Value = 5138.0769
$ws.Range("M122").Value = -2688.0769
$ws.Range("H126").Value = 814.7826
$ws.Range("I126").Value = 701.9048
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 2105.7144
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = 364.2856000000002
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2647.36
$ws.Range("I132").Value = 2239.25
$ws.Range("J132").Value = 3024.077
$ws.Range("K132").Value = 6717.75
$ws.Range("L132").Value = 9072.231
$ws.Range("M132").Value = -4187.75
$ws.Range("N132").Value = -14132.231
$ws.Range("H136").Value = 2709.5103
$ws.Range("I136").Value = 3054.2083
$ws.Range("K136").Value = 9162.624899999999
$ws.Range("M136").Value = -6612.624899999999
